# Adapt column header formatting to respective input file names:
#   *_old -> *_FV2404
#   *_new -> *_FV2410
# and turn the header row (A1:U1) into an Excel Table (ListObject) with a
# frozen header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header names for columns A..U (in order).
$headers = @(
  "Segmentname_FV2404",
  "Segmentgruppe_FV2404",
  "Segment_FV2404",
  "Datenelement_FV2404",
  "Segment ID_FV2404",
  "Code_FV2404",
  "Qualifier_FV2404",
  "Beschreibung_FV2404",
  "Bedingungsausdruck_FV2404",
  "Bedingung_FV2404",
  "diff",
  "Segmentname_FV2410",
  "Segmentgruppe_FV2410",
  "Segment_FV2410",
  "Datenelement_FV2410",
  "Segment ID_FV2410",
  "Code_FV2410",
  "Qualifier_FV2410",
  "Beschreibung_FV2410",
  "Bedingungsausdruck_FV2410",
  "Bedingung_FV2410"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Stash the existing header formatting (bold, fill, border, alignment) into
# unused scratch cells on the same sheet so it survives table creation
# unchanged (creating a ListObject over an already-formatted header row
# would otherwise register a conditional-formatting "header row" dxf).
$headerRange = $ws.Range("A1:U1")
$scratchRange = $ws.Range("AA1:AU1")

$headerRange.Copy()
$scratchRange.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$headerRange.ClearFormats()

# Turn A1:U78 into a proper Excel Table.
$tableRange = $ws.Range("A1:U78")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)  # xlSrcRange, xlYes
$tbl.TableStyle = ""

# Restore the original header formatting onto the (now table) header row.
$scratchRange.Copy()
$headerRange.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Remove the scratch cells again.
$scratchRange.Clear()

# Freeze the header row (split below row 1).
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
